# "Transport von bis zu 4 Personen eingefügt"
# Adds 8 more order rows (rows 10-17) to the "Orders" sheet (duplicating the
# existing order in rows 2-9, i.e. up to 4 taxi requests x 2 passengers each)
# and a new taxi (row 3) to the "Taxis" sheet, then makes "Taxis" the active tab.

$wb = $excel.ActiveWorkbook

# --- Orders sheet: duplicate the order block into rows 10-17 -------------
$wsOrders = $wb.Worksheets.Item("Orders")

for ($r = 10; $r -le 17; $r++) {
    $wsOrders.Cells.Item($r, 1).Value = 120
    $wsOrders.Cells.Item($r, 2).Value = 52.375394479042797
    $wsOrders.Cells.Item($r, 3).Value = 9.7315180260351593
    $wsOrders.Cells.Item($r, 4).Value = 52.382591097574597
    $wsOrders.Cells.Item($r, 5).Value = 9.7309718027690906
}

$wsOrders.Range("B9:B17").Select()
$wsOrders.Application.ActiveWindow.RangeSelection.Item(1).Activate()

# --- TaxiStands sheet: move selection (no data change) -------------------
$wsStands = $wb.Worksheets.Item("TaxiStands")
$wsStands.Range("A3:B3").Select()
$wsStands.Application.ActiveWindow.RangeSelection.Item(1).Activate()

# --- Taxis sheet: add a new taxi in row 3 ---------------------------------
$wsTaxis = $wb.Worksheets.Item("Taxis")

$wsTaxis.Cells.Item(3, 1).Value = 52.389690000000002
$wsTaxis.Cells.Item(3, 2).Value = 9.7207399999999993
$wsTaxis.Cells.Item(3, 4).Value = 30000
$wsTaxis.Cells.Item(3, 5).Value = 30000
$wsTaxis.Cells.Item(3, 6).Value = 30000
$wsTaxis.Cells.Item(3, 7).Value = "Tesla"

$wsTaxis.Activate()
$wsTaxis.Range("G3").Select()

Write-Host "done"
